# Adds the Frequency-mode display-measure translation rows (TEXT ID, TYPOGRAPHY NAME,
# ALIGNMENT, DIRECTION, GB) to the "Translation" sheet, rows 184-273, as introduced by
# the "Developed UI for display measures Frequency Mode." commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 184
$ws.Cells.Item(184, 2).Value = "SingleUseId199"
$ws.Cells.Item(184, 3).Value = "displayMeas"
$ws.Cells.Item(184, 4).Value = "Left"
$ws.Cells.Item(184, 5).Value = "LTR"
$ws.Cells.Item(184, 6).Value = "(IN<value>)"

# Row 185
$ws.Cells.Item(185, 2).Value = "SingleUseId200"
$ws.Cells.Item(185, 3).Value = "displayMeas"
$ws.Cells.Item(185, 4).Value = "Left"
$ws.Cells.Item(185, 5).Value = "LTR"
$ws.Cells.Item(185, 6).Value = "(IN<value>)"

# Row 186
$ws.Cells.Item(186, 2).Value = "SingleUseId201"
$ws.Cells.Item(186, 3).Value = "displayMeas"
$ws.Cells.Item(186, 4).Value = "Left"
$ws.Cells.Item(186, 5).Value = "LTR"
$ws.Cells.Item(186, 6).Value = "(IN<value>)"

# Row 187
$ws.Cells.Item(187, 2).Value = "SingleUseId202"
$ws.Cells.Item(187, 3).Value = "displayMeas"
$ws.Cells.Item(187, 4).Value = "Left"
$ws.Cells.Item(187, 5).Value = "LTR"
$ws.Cells.Item(187, 6).Value = "(IN<value>)"

# Row 188
$ws.Cells.Item(188, 2).Value = "SingleUseId203"
$ws.Cells.Item(188, 3).Value = "displayMeas"
$ws.Cells.Item(188, 4).Value = "Left"
$ws.Cells.Item(188, 5).Value = "LTR"
$ws.Cells.Item(188, 6).Value = "(IN<value>)"

# Row 189
$ws.Cells.Item(189, 2).Value = "SingleUseId204"
$ws.Cells.Item(189, 3).Value = "displayMeas"
$ws.Cells.Item(189, 4).Value = "Left"
$ws.Cells.Item(189, 5).Value = "LTR"
$ws.Cells.Item(189, 6).Value = "(IN<value>)"

# Row 190
$ws.Cells.Item(190, 2).Value = "SingleUseId205"
$ws.Cells.Item(190, 3).Value = "displayMeas"
$ws.Cells.Item(190, 4).Value = "Left"
$ws.Cells.Item(190, 5).Value = "LTR"
$ws.Cells.Item(190, 6).Value = "(IN<value>)"

# Row 191
$ws.Cells.Item(191, 2).Value = "SingleUseId206"
$ws.Cells.Item(191, 3).Value = "displayMeas"
$ws.Cells.Item(191, 4).Value = "Left"
$ws.Cells.Item(191, 5).Value = "LTR"
$ws.Cells.Item(191, 6).Value = "Freq 1 "

# Row 192
$ws.Cells.Item(192, 2).Value = "SingleUseId207"
$ws.Cells.Item(192, 3).Value = "displayMeas"
$ws.Cells.Item(192, 4).Value = "Left"
$ws.Cells.Item(192, 5).Value = "LTR"
$ws.Cells.Item(192, 6).Value = "Freq 2"

# Row 193
$ws.Cells.Item(193, 2).Value = "SingleUseId208"
$ws.Cells.Item(193, 3).Value = "displayMeas"
$ws.Cells.Item(193, 4).Value = "Left"
$ws.Cells.Item(193, 5).Value = "LTR"
$ws.Cells.Item(193, 6).Value = "Freq 3"

# Row 194
$ws.Cells.Item(194, 2).Value = "SingleUseId209"
$ws.Cells.Item(194, 3).Value = "displayMeas"
$ws.Cells.Item(194, 4).Value = "Left"
$ws.Cells.Item(194, 5).Value = "LTR"
$ws.Cells.Item(194, 6).Value = "Freq 4"

# Row 195
$ws.Cells.Item(195, 2).Value = "SingleUseId210"
$ws.Cells.Item(195, 3).Value = "displayMeas"
$ws.Cells.Item(195, 4).Value = "Left"
$ws.Cells.Item(195, 5).Value = "LTR"
$ws.Cells.Item(195, 6).Value = "Freq 5"

# Row 196
$ws.Cells.Item(196, 2).Value = "SingleUseId211"
$ws.Cells.Item(196, 3).Value = "displayMeas"
$ws.Cells.Item(196, 4).Value = "Left"
$ws.Cells.Item(196, 5).Value = "LTR"
$ws.Cells.Item(196, 6).Value = "Freq 6"

# Row 197
$ws.Cells.Item(197, 2).Value = "SingleUseId212"
$ws.Cells.Item(197, 3).Value = "displayMeas"
$ws.Cells.Item(197, 4).Value = "Left"
$ws.Cells.Item(197, 5).Value = "LTR"
$ws.Cells.Item(197, 6).Value = "Freq 7"

# Row 198
$ws.Cells.Item(198, 2).Value = "SingleUseId213"
$ws.Cells.Item(198, 3).Value = "displayMeas"
$ws.Cells.Item(198, 4).Value = "Left"
$ws.Cells.Item(198, 5).Value = "LTR"
$ws.Cells.Item(198, 6).Value = "(IN<value>)"

# Row 199
$ws.Cells.Item(199, 2).Value = "SingleUseId214"
$ws.Cells.Item(199, 3).Value = "displayMeas"
$ws.Cells.Item(199, 4).Value = "Left"
$ws.Cells.Item(199, 5).Value = "LTR"
$ws.Cells.Item(199, 6).Value = "Freq 8"

# Row 200
$ws.Cells.Item(200, 2).Value = "SingleUseId215"
$ws.Cells.Item(200, 3).Value = "displaylabelBold"
$ws.Cells.Item(200, 4).Value = "Left"
$ws.Cells.Item(200, 5).Value = "LTR"
$ws.Cells.Item(200, 6).Value = "Sample"

# Row 201
$ws.Cells.Item(201, 2).Value = "SingleUseId216"
$ws.Cells.Item(201, 3).Value = "displayMeas"
$ws.Cells.Item(201, 4).Value = "Left"
$ws.Cells.Item(201, 5).Value = "LTR"
$ws.Cells.Item(201, 6).Value = "Mean"

# Row 202
$ws.Cells.Item(202, 2).Value = "SingleUseId217"
$ws.Cells.Item(202, 3).Value = "displayMeas"
$ws.Cells.Item(202, 4).Value = "Left"
$ws.Cells.Item(202, 5).Value = "LTR"
$ws.Cells.Item(202, 6).Value = "StdDev"

# Row 203
$ws.Cells.Item(203, 2).Value = "SingleUseId218"
$ws.Cells.Item(203, 3).Value = "displayLabel"
$ws.Cells.Item(203, 4).Value = "Center"
$ws.Cells.Item(203, 5).Value = "LTR"
$ws.Cells.Item(203, 6).Value = "<value>"

# Row 204
$ws.Cells.Item(204, 2).Value = "SingleUseId219"
$ws.Cells.Item(204, 3).Value = "displayLabel"
$ws.Cells.Item(204, 4).Value = "Center"
$ws.Cells.Item(204, 5).Value = "LTR"
$ws.Cells.Item(204, 6).Value = "<value>"

# Row 205
$ws.Cells.Item(205, 2).Value = "SingleUseId220"
$ws.Cells.Item(205, 3).Value = "displayLabel"
$ws.Cells.Item(205, 4).Value = "Center"
$ws.Cells.Item(205, 5).Value = "LTR"
$ws.Cells.Item(205, 6).Value = "<value>"

# Row 206
$ws.Cells.Item(206, 2).Value = "SingleUseId221"
$ws.Cells.Item(206, 3).Value = "displayLabel"
$ws.Cells.Item(206, 4).Value = "Center"
$ws.Cells.Item(206, 5).Value = "LTR"
$ws.Cells.Item(206, 6).Value = "<value>"

# Row 207
$ws.Cells.Item(207, 2).Value = "SingleUseId222"
$ws.Cells.Item(207, 3).Value = "displayLabel"
$ws.Cells.Item(207, 4).Value = "Center"
$ws.Cells.Item(207, 5).Value = "LTR"
$ws.Cells.Item(207, 6).Value = "<value>"

# Row 208
$ws.Cells.Item(208, 2).Value = "SingleUseId223"
$ws.Cells.Item(208, 3).Value = "displayLabel"
$ws.Cells.Item(208, 4).Value = "Center"
$ws.Cells.Item(208, 5).Value = "LTR"
$ws.Cells.Item(208, 6).Value = "<value>"

# Row 209
$ws.Cells.Item(209, 2).Value = "SingleUseId224"
$ws.Cells.Item(209, 3).Value = "displayLabel"
$ws.Cells.Item(209, 4).Value = "Center"
$ws.Cells.Item(209, 5).Value = "LTR"
$ws.Cells.Item(209, 6).Value = "<value>"

# Row 210
$ws.Cells.Item(210, 2).Value = "SingleUseId225"
$ws.Cells.Item(210, 3).Value = "displayLabel"
$ws.Cells.Item(210, 4).Value = "Center"
$ws.Cells.Item(210, 5).Value = "LTR"
$ws.Cells.Item(210, 6).Value = "<value>"

# Row 211
$ws.Cells.Item(211, 2).Value = "SingleUseId226"
$ws.Cells.Item(211, 3).Value = "displayMeasValue"
$ws.Cells.Item(211, 4).Value = "Right"
$ws.Cells.Item(211, 5).Value = "LTR"
$ws.Cells.Item(211, 6).Value = "<value>"

# Row 212
$ws.Cells.Item(212, 2).Value = "SingleUseId227"
$ws.Cells.Item(212, 3).Value = "displayMeasValue"
$ws.Cells.Item(212, 4).Value = "Right"
$ws.Cells.Item(212, 5).Value = "LTR"
$ws.Cells.Item(212, 6).Value = "<value>"

# Row 213
$ws.Cells.Item(213, 2).Value = "SingleUseId228"
$ws.Cells.Item(213, 3).Value = "displayMeasValue"
$ws.Cells.Item(213, 4).Value = "Right"
$ws.Cells.Item(213, 5).Value = "LTR"
$ws.Cells.Item(213, 6).Value = "<value>"

# Row 214
$ws.Cells.Item(214, 2).Value = "SingleUseId229"
$ws.Cells.Item(214, 3).Value = "displayMeasValue"
$ws.Cells.Item(214, 4).Value = "Right"
$ws.Cells.Item(214, 5).Value = "LTR"
$ws.Cells.Item(214, 6).Value = "<value>"

# Row 215
$ws.Cells.Item(215, 2).Value = "SingleUseId230"
$ws.Cells.Item(215, 3).Value = "displayMeasValue"
$ws.Cells.Item(215, 4).Value = "Right"
$ws.Cells.Item(215, 5).Value = "LTR"
$ws.Cells.Item(215, 6).Value = "<value>"

# Row 216
$ws.Cells.Item(216, 2).Value = "SingleUseId231"
$ws.Cells.Item(216, 3).Value = "displayMeasValue"
$ws.Cells.Item(216, 4).Value = "Right"
$ws.Cells.Item(216, 5).Value = "LTR"
$ws.Cells.Item(216, 6).Value = "<value>"

# Row 217
$ws.Cells.Item(217, 2).Value = "SingleUseId232"
$ws.Cells.Item(217, 3).Value = "displayMeasValue"
$ws.Cells.Item(217, 4).Value = "Right"
$ws.Cells.Item(217, 5).Value = "LTR"
$ws.Cells.Item(217, 6).Value = "<value>"

# Row 218
$ws.Cells.Item(218, 2).Value = "SingleUseId233"
$ws.Cells.Item(218, 3).Value = "displayMeasValue"
$ws.Cells.Item(218, 4).Value = "Right"
$ws.Cells.Item(218, 5).Value = "LTR"
$ws.Cells.Item(218, 6).Value = "<value>"

# Row 219
$ws.Cells.Item(219, 2).Value = "SingleUseId234"
$ws.Cells.Item(219, 3).Value = "displayMeasValue"
$ws.Cells.Item(219, 4).Value = "Right"
$ws.Cells.Item(219, 5).Value = "LTR"
$ws.Cells.Item(219, 6).Value = "<value>"

# Row 220
$ws.Cells.Item(220, 2).Value = "SingleUseId235"
$ws.Cells.Item(220, 3).Value = "displayMeasValue"
$ws.Cells.Item(220, 4).Value = "Right"
$ws.Cells.Item(220, 5).Value = "LTR"
$ws.Cells.Item(220, 6).Value = "<value>"

# Row 221
$ws.Cells.Item(221, 2).Value = "SingleUseId236"
$ws.Cells.Item(221, 3).Value = "displayMeasValue"
$ws.Cells.Item(221, 4).Value = "Right"
$ws.Cells.Item(221, 5).Value = "LTR"
$ws.Cells.Item(221, 6).Value = "<value>"

# Row 222
$ws.Cells.Item(222, 2).Value = "SingleUseId237"
$ws.Cells.Item(222, 3).Value = "displayMeasValue"
$ws.Cells.Item(222, 4).Value = "Right"
$ws.Cells.Item(222, 5).Value = "LTR"
$ws.Cells.Item(222, 6).Value = "<value>"

# Row 223
$ws.Cells.Item(223, 2).Value = "SingleUseId238"
$ws.Cells.Item(223, 3).Value = "displayMeasValue"
$ws.Cells.Item(223, 4).Value = "Right"
$ws.Cells.Item(223, 5).Value = "LTR"
$ws.Cells.Item(223, 6).Value = "<value>"

# Row 224
$ws.Cells.Item(224, 2).Value = "SingleUseId239"
$ws.Cells.Item(224, 3).Value = "displayMeasValue"
$ws.Cells.Item(224, 4).Value = "Right"
$ws.Cells.Item(224, 5).Value = "LTR"
$ws.Cells.Item(224, 6).Value = "<value>"

# Row 225
$ws.Cells.Item(225, 2).Value = "SingleUseId240"
$ws.Cells.Item(225, 3).Value = "displayMeasValue"
$ws.Cells.Item(225, 4).Value = "Right"
$ws.Cells.Item(225, 5).Value = "LTR"
$ws.Cells.Item(225, 6).Value = "<value>"

# Row 226
$ws.Cells.Item(226, 2).Value = "SingleUseId241"
$ws.Cells.Item(226, 3).Value = "displayMeasValue"
$ws.Cells.Item(226, 4).Value = "Right"
$ws.Cells.Item(226, 5).Value = "LTR"
$ws.Cells.Item(226, 6).Value = "<value>"

# Row 227
$ws.Cells.Item(227, 2).Value = "SingleUseId242"
$ws.Cells.Item(227, 3).Value = "displayLabel"
$ws.Cells.Item(227, 4).Value = "Left"
$ws.Cells.Item(227, 5).Value = "LTR"
$ws.Cells.Item(227, 6).Value = "ns"

# Row 228
$ws.Cells.Item(228, 2).Value = "SingleUseId243"
$ws.Cells.Item(228, 3).Value = "displayLabel"
$ws.Cells.Item(228, 4).Value = "Left"
$ws.Cells.Item(228, 5).Value = "LTR"
$ws.Cells.Item(228, 6).Value = "ns"

# Row 229
$ws.Cells.Item(229, 2).Value = "SingleUseId244"
$ws.Cells.Item(229, 3).Value = "displayLabel"
$ws.Cells.Item(229, 4).Value = "Left"
$ws.Cells.Item(229, 5).Value = "LTR"
$ws.Cells.Item(229, 6).Value = "ns"

# Row 230
$ws.Cells.Item(230, 2).Value = "SingleUseId245"
$ws.Cells.Item(230, 3).Value = "displayLabel"
$ws.Cells.Item(230, 4).Value = "Left"
$ws.Cells.Item(230, 5).Value = "LTR"
$ws.Cells.Item(230, 6).Value = "ns"

# Row 231
$ws.Cells.Item(231, 2).Value = "SingleUseId246"
$ws.Cells.Item(231, 3).Value = "displayLabel"
$ws.Cells.Item(231, 4).Value = "Left"
$ws.Cells.Item(231, 5).Value = "LTR"
$ws.Cells.Item(231, 6).Value = "ns"

# Row 232
$ws.Cells.Item(232, 2).Value = "SingleUseId247"
$ws.Cells.Item(232, 3).Value = "displayLabel"
$ws.Cells.Item(232, 4).Value = "Left"
$ws.Cells.Item(232, 5).Value = "LTR"
$ws.Cells.Item(232, 6).Value = "ns"

# Row 233
$ws.Cells.Item(233, 2).Value = "SingleUseId248"
$ws.Cells.Item(233, 3).Value = "displayLabel"
$ws.Cells.Item(233, 4).Value = "Left"
$ws.Cells.Item(233, 5).Value = "LTR"
$ws.Cells.Item(233, 6).Value = "ns"

# Row 234
$ws.Cells.Item(234, 2).Value = "SingleUseId249"
$ws.Cells.Item(234, 3).Value = "displayLabel"
$ws.Cells.Item(234, 4).Value = "Left"
$ws.Cells.Item(234, 5).Value = "LTR"
$ws.Cells.Item(234, 6).Value = "uHz"

# Row 235
$ws.Cells.Item(235, 2).Value = "SingleUseId250"
$ws.Cells.Item(235, 3).Value = "displayLabel"
$ws.Cells.Item(235, 4).Value = "Left"
$ws.Cells.Item(235, 5).Value = "LTR"
$ws.Cells.Item(235, 6).Value = "uHz"

# Row 236
$ws.Cells.Item(236, 2).Value = "SingleUseId251"
$ws.Cells.Item(236, 3).Value = "displayLabel"
$ws.Cells.Item(236, 4).Value = "Left"
$ws.Cells.Item(236, 5).Value = "LTR"
$ws.Cells.Item(236, 6).Value = "uHz"

# Row 237
$ws.Cells.Item(237, 2).Value = "SingleUseId252"
$ws.Cells.Item(237, 3).Value = "displayLabel"
$ws.Cells.Item(237, 4).Value = "Left"
$ws.Cells.Item(237, 5).Value = "LTR"
$ws.Cells.Item(237, 6).Value = "uHz"

# Row 238
$ws.Cells.Item(238, 2).Value = "SingleUseId253"
$ws.Cells.Item(238, 3).Value = "displayLabel"
$ws.Cells.Item(238, 4).Value = "Left"
$ws.Cells.Item(238, 5).Value = "LTR"
$ws.Cells.Item(238, 6).Value = "uHz"

# Row 239
$ws.Cells.Item(239, 2).Value = "SingleUseId254"
$ws.Cells.Item(239, 3).Value = "displayLabel"
$ws.Cells.Item(239, 4).Value = "Left"
$ws.Cells.Item(239, 5).Value = "LTR"
$ws.Cells.Item(239, 6).Value = "uHz"

# Row 240
$ws.Cells.Item(240, 2).Value = "SingleUseId255"
$ws.Cells.Item(240, 3).Value = "displayLabel"
$ws.Cells.Item(240, 4).Value = "Left"
$ws.Cells.Item(240, 5).Value = "LTR"
$ws.Cells.Item(240, 6).Value = "uHz"

# Row 241
$ws.Cells.Item(241, 2).Value = "SingleUseId256"
$ws.Cells.Item(241, 3).Value = "displayLabel"
$ws.Cells.Item(241, 4).Value = "Left"
$ws.Cells.Item(241, 5).Value = "LTR"
$ws.Cells.Item(241, 6).Value = "uHz"

# Row 242
$ws.Cells.Item(242, 2).Value = "SingleUseId257"
$ws.Cells.Item(242, 3).Value = "displayLabel"
$ws.Cells.Item(242, 4).Value = "Left"
$ws.Cells.Item(242, 5).Value = "LTR"
$ws.Cells.Item(242, 6).Value = "mHz"

# Row 243
$ws.Cells.Item(243, 2).Value = "SingleUseId258"
$ws.Cells.Item(243, 3).Value = "displayLabel"
$ws.Cells.Item(243, 4).Value = "Left"
$ws.Cells.Item(243, 5).Value = "LTR"
$ws.Cells.Item(243, 6).Value = "mHz"

# Row 244
$ws.Cells.Item(244, 2).Value = "SingleUseId259"
$ws.Cells.Item(244, 3).Value = "displayLabel"
$ws.Cells.Item(244, 4).Value = "Left"
$ws.Cells.Item(244, 5).Value = "LTR"
$ws.Cells.Item(244, 6).Value = "mHz"

# Row 245
$ws.Cells.Item(245, 2).Value = "SingleUseId260"
$ws.Cells.Item(245, 3).Value = "displayLabel"
$ws.Cells.Item(245, 4).Value = "Left"
$ws.Cells.Item(245, 5).Value = "LTR"
$ws.Cells.Item(245, 6).Value = "mHz"

# Row 246
$ws.Cells.Item(246, 2).Value = "SingleUseId261"
$ws.Cells.Item(246, 3).Value = "displayLabel"
$ws.Cells.Item(246, 4).Value = "Left"
$ws.Cells.Item(246, 5).Value = "LTR"
$ws.Cells.Item(246, 6).Value = "mHz"

# Row 247
$ws.Cells.Item(247, 2).Value = "SingleUseId262"
$ws.Cells.Item(247, 3).Value = "displayLabel"
$ws.Cells.Item(247, 4).Value = "Left"
$ws.Cells.Item(247, 5).Value = "LTR"
$ws.Cells.Item(247, 6).Value = "mHz"

# Row 248
$ws.Cells.Item(248, 2).Value = "SingleUseId263"
$ws.Cells.Item(248, 3).Value = "displayLabel"
$ws.Cells.Item(248, 4).Value = "Left"
$ws.Cells.Item(248, 5).Value = "LTR"
$ws.Cells.Item(248, 6).Value = "mHz"

# Row 249
$ws.Cells.Item(249, 2).Value = "SingleUseId264"
$ws.Cells.Item(249, 3).Value = "displayLabel"
$ws.Cells.Item(249, 4).Value = "Left"
$ws.Cells.Item(249, 5).Value = "LTR"
$ws.Cells.Item(249, 6).Value = "mHz"

# Row 250
$ws.Cells.Item(250, 2).Value = "SingleUseId265"
$ws.Cells.Item(250, 3).Value = "displayLabel"
$ws.Cells.Item(250, 4).Value = "Left"
$ws.Cells.Item(250, 5).Value = "LTR"
$ws.Cells.Item(250, 6).Value = "kHz"

# Row 251
$ws.Cells.Item(251, 2).Value = "SingleUseId266"
$ws.Cells.Item(251, 3).Value = "displayLabel"
$ws.Cells.Item(251, 4).Value = "Left"
$ws.Cells.Item(251, 5).Value = "LTR"
$ws.Cells.Item(251, 6).Value = "kHz"

# Row 252
$ws.Cells.Item(252, 2).Value = "SingleUseId267"
$ws.Cells.Item(252, 3).Value = "displayLabel"
$ws.Cells.Item(252, 4).Value = "Left"
$ws.Cells.Item(252, 5).Value = "LTR"
$ws.Cells.Item(252, 6).Value = "kHz"

# Row 253
$ws.Cells.Item(253, 2).Value = "SingleUseId268"
$ws.Cells.Item(253, 3).Value = "displayLabel"
$ws.Cells.Item(253, 4).Value = "Left"
$ws.Cells.Item(253, 5).Value = "LTR"
$ws.Cells.Item(253, 6).Value = "kHz"

# Row 254
$ws.Cells.Item(254, 2).Value = "SingleUseId269"
$ws.Cells.Item(254, 3).Value = "displayLabel"
$ws.Cells.Item(254, 4).Value = "Left"
$ws.Cells.Item(254, 5).Value = "LTR"
$ws.Cells.Item(254, 6).Value = "kHz"

# Row 255
$ws.Cells.Item(255, 2).Value = "SingleUseId270"
$ws.Cells.Item(255, 3).Value = "displayLabel"
$ws.Cells.Item(255, 4).Value = "Left"
$ws.Cells.Item(255, 5).Value = "LTR"
$ws.Cells.Item(255, 6).Value = "kHz"

# Row 256
$ws.Cells.Item(256, 2).Value = "SingleUseId271"
$ws.Cells.Item(256, 3).Value = "displayLabel"
$ws.Cells.Item(256, 4).Value = "Left"
$ws.Cells.Item(256, 5).Value = "LTR"
$ws.Cells.Item(256, 6).Value = "kHz"

# Row 257
$ws.Cells.Item(257, 2).Value = "SingleUseId272"
$ws.Cells.Item(257, 3).Value = "displayLabel"
$ws.Cells.Item(257, 4).Value = "Left"
$ws.Cells.Item(257, 5).Value = "LTR"
$ws.Cells.Item(257, 6).Value = "kHz"

# Row 258
$ws.Cells.Item(258, 2).Value = "SingleUseId273"
$ws.Cells.Item(258, 3).Value = "displayLabel"
$ws.Cells.Item(258, 4).Value = "Left"
$ws.Cells.Item(258, 5).Value = "LTR"
$ws.Cells.Item(258, 6).Value = "MHz"

# Row 259
$ws.Cells.Item(259, 2).Value = "SingleUseId274"
$ws.Cells.Item(259, 3).Value = "displayLabel"
$ws.Cells.Item(259, 4).Value = "Left"
$ws.Cells.Item(259, 5).Value = "LTR"
$ws.Cells.Item(259, 6).Value = "MHz"

# Row 260
$ws.Cells.Item(260, 2).Value = "SingleUseId275"
$ws.Cells.Item(260, 3).Value = "displayLabel"
$ws.Cells.Item(260, 4).Value = "Left"
$ws.Cells.Item(260, 5).Value = "LTR"
$ws.Cells.Item(260, 6).Value = "MHz"

# Row 261
$ws.Cells.Item(261, 2).Value = "SingleUseId276"
$ws.Cells.Item(261, 3).Value = "displayLabel"
$ws.Cells.Item(261, 4).Value = "Left"
$ws.Cells.Item(261, 5).Value = "LTR"
$ws.Cells.Item(261, 6).Value = "MHz"

# Row 262
$ws.Cells.Item(262, 2).Value = "SingleUseId277"
$ws.Cells.Item(262, 3).Value = "displayLabel"
$ws.Cells.Item(262, 4).Value = "Left"
$ws.Cells.Item(262, 5).Value = "LTR"
$ws.Cells.Item(262, 6).Value = "MHz"

# Row 263
$ws.Cells.Item(263, 2).Value = "SingleUseId278"
$ws.Cells.Item(263, 3).Value = "displayLabel"
$ws.Cells.Item(263, 4).Value = "Left"
$ws.Cells.Item(263, 5).Value = "LTR"
$ws.Cells.Item(263, 6).Value = "MHz"

# Row 264
$ws.Cells.Item(264, 2).Value = "SingleUseId279"
$ws.Cells.Item(264, 3).Value = "displayLabel"
$ws.Cells.Item(264, 4).Value = "Left"
$ws.Cells.Item(264, 5).Value = "LTR"
$ws.Cells.Item(264, 6).Value = "MHz"

# Row 265
$ws.Cells.Item(265, 2).Value = "SingleUseId280"
$ws.Cells.Item(265, 3).Value = "displayLabel"
$ws.Cells.Item(265, 4).Value = "Left"
$ws.Cells.Item(265, 5).Value = "LTR"
$ws.Cells.Item(265, 6).Value = "MHz"

# Row 266
$ws.Cells.Item(266, 2).Value = "SingleUseId281"
$ws.Cells.Item(266, 3).Value = "displayLabel"
$ws.Cells.Item(266, 4).Value = "Left"
$ws.Cells.Item(266, 5).Value = "LTR"
$ws.Cells.Item(266, 6).Value = "GHz"

# Row 267
$ws.Cells.Item(267, 2).Value = "SingleUseId282"
$ws.Cells.Item(267, 3).Value = "displayLabel"
$ws.Cells.Item(267, 4).Value = "Left"
$ws.Cells.Item(267, 5).Value = "LTR"
$ws.Cells.Item(267, 6).Value = "GHz"

# Row 268
$ws.Cells.Item(268, 2).Value = "SingleUseId283"
$ws.Cells.Item(268, 3).Value = "displayLabel"
$ws.Cells.Item(268, 4).Value = "Left"
$ws.Cells.Item(268, 5).Value = "LTR"
$ws.Cells.Item(268, 6).Value = "GHz"

# Row 269
$ws.Cells.Item(269, 2).Value = "SingleUseId284"
$ws.Cells.Item(269, 3).Value = "displayLabel"
$ws.Cells.Item(269, 4).Value = "Left"
$ws.Cells.Item(269, 5).Value = "LTR"
$ws.Cells.Item(269, 6).Value = "GHz"

# Row 270
$ws.Cells.Item(270, 2).Value = "SingleUseId285"
$ws.Cells.Item(270, 3).Value = "displayLabel"
$ws.Cells.Item(270, 4).Value = "Left"
$ws.Cells.Item(270, 5).Value = "LTR"
$ws.Cells.Item(270, 6).Value = "GHz"

# Row 271
$ws.Cells.Item(271, 2).Value = "SingleUseId286"
$ws.Cells.Item(271, 3).Value = "displayLabel"
$ws.Cells.Item(271, 4).Value = "Left"
$ws.Cells.Item(271, 5).Value = "LTR"
$ws.Cells.Item(271, 6).Value = "GHz"

# Row 272
$ws.Cells.Item(272, 2).Value = "SingleUseId287"
$ws.Cells.Item(272, 3).Value = "displayLabel"
$ws.Cells.Item(272, 4).Value = "Left"
$ws.Cells.Item(272, 5).Value = "LTR"
$ws.Cells.Item(272, 6).Value = "GHz"

# Row 273
$ws.Cells.Item(273, 2).Value = "SingleUseId288"
$ws.Cells.Item(273, 3).Value = "displayLabel"
$ws.Cells.Item(273, 4).Value = "Left"
$ws.Cells.Item(273, 5).Value = "LTR"
$ws.Cells.Item(273, 6).Value = "GHz"
